$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.638.00'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.763.66'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.21'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3825'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3421'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.57'
$ws.Range('E9').Value = '  -3.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.133'
$ws.Range('E10').Value = '  -5.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07389'
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.23'
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.333'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.764.45'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.063'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001072'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06653'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '81.88'
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.0000'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.28'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.386'
$ws.Range('E22').Value = '  -3.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.644.37'
$ws.Range('E23').Value = '  +1.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.99'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.388'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.424'
$ws.Range('E26').Value = '  -3.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.58'
$ws.Range('E27').Value = '  -2.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.401'
$ws.Range('E28').Value = '  -5.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '152.53'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.964.81'
$ws.Range('E30').Value = '  -1.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '133.65'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.083'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.954'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08774'
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.67'
$ws.Range('E35').Value = '  -4.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02410'
$ws.Range('E36').Value = '  +3.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6764'
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.295'
$ws.Range('E38').Value = '  -2.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06284'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2174'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.252'
$ws.Range('E41').Value = '  +1.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.499'
$ws.Range('E42').Value = '  -8.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.202'
$ws.Range('E43').Value = '  -6.31%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.04'
$ws.Range('E45').Value = '  -1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6236'
$ws.Range('E46').Value = '  -3.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.837'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.23'
$ws.Range('E48').Value = '  +0.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.065'
$ws.Range('E49').Value = '  -3.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07401'
$ws.Range('E50').Value = '  +3.93%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.141'
$ws.Range('E51').Value = '  +2.23%  '
